$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 342
$ws.Range("I4").Value = 342
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 342
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -228
$ws.Range("N4").ClearContents()
$ws.Range("H6").Value = 433.9
$ws.Range("I6").Value = 218
$ws.Range("J6").Value = 757.75
$ws.Range("K6").Value = 654
$ws.Range("L6").Value = 2273.25
$ws.Range("M6").Value = -542
$ws.Range("N6").Value = -2497.25
$ws.Range("H20").Value = 550
$ws.Range("I20").Value = 550
$ws.Range("K20").Value = 550
$ws.Range("M20").Value = -320
$ws.Range("H32").Value = 1038
$ws.Range("I32").Value = 813.6667
$ws.Range("J32").Value = 1374.5
$ws.Range("K32").Value = 813.6667
$ws.Range("L32").Value = 1374.5
$ws.Range("M32").Value = -487.6667
$ws.Range("N32").Value = -2026.5
$ws.Range("H35").Value = 550
$ws.Range("I35").Value = 550
$ws.Range("K35").Value = 550
$ws.Range("M35").Value = -171
$ws.Range("H39").Value = 3959
$ws.Range("I39").Value = 23.75
$ws.Range("J39").Value = 19700
$ws.Range("K39").Value = 71.25
$ws.Range("L39").Value = 59100
$ws.Range("M39").Value = 224.75
$ws.Range("N39").Value = -59692
$ws.Range("H42").Value = 881.6667
$ws.Range("I42").Value = 63.333332
$ws.Range("K42").Value = 189.999996
$ws.Range("M42").Value = 40.00000399999999
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H58").Value = 1341.5
$ws.Range("I58").Value = 1303.75
$ws.Range("J58").Value = 1366.6666
$ws.Range("K58").Value = 3911.25
$ws.Range("L58").Value = 4099.9998
$ws.Range("M58").Value = -3761.25
$ws.Range("N58").Value = -4399.9998
$ws.Range("H98").Value = 3850.8572
$ws.Range("I98").Value = 3660.5
$ws.Range("J98").Value = 4993
$ws.Range("K98").Value = 3660.5
$ws.Range("L98").Value = 4993
$ws.Range("M98").Value = -2162.5
$ws.Range("N98").Value = -7989
$ws.Range("H100").Value = 1500
$ws.Range("I100").Value = 1000
$ws.Range("K100").Value = 1000
$ws.Range("M100").Value = -459
$ws.Range("H122").Value = 3850.8572
$ws.Range("I122").Value = 3660.5
$ws.Range("J122").Value = 4993
$ws.Range("K122").Value = 10981.5
$ws.Range("L122").Value = 14979
$ws.Range("M122").Value = -8531.5
$ws.Range("N122").Value = -19879

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5000131.5
$ws.Range("I6").Value = 5000131.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 5000131.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -5000018.5
$ws.Range("N6").ClearContents()
$ws.Range("H19").Value = 279.82352
$ws.Range("I19").Value = 165.11111
$ws.Range("K19").Value = 165.11111
$ws.Range("M19").Value = 4.888890000000004
$ws.Range("H24").Value = 279.82352
$ws.Range("I24").Value = 165.11111
$ws.Range("K24").Value = 165.11111
$ws.Range("M24").Value = 4.888890000000004
$ws.Range("H39").Value = 3269
$ws.Range("I39").Value = 3375
$ws.Range("J39").Value = 3057
$ws.Range("K39").Value = 3375
$ws.Range("L39").Value = 3057
$ws.Range("M39").Value = -2984
$ws.Range("N39").Value = -3839
$ws.Range("H49").Value = 3269
$ws.Range("I49").Value = 3375
$ws.Range("J49").Value = 3057
$ws.Range("K49").Value = 3375
$ws.Range("L49").Value = 3057
$ws.Range("M49").Value = -3193
$ws.Range("N49").Value = -3421
$ws.Range("H74").Value = 47500
$ws.Range("J74").Value = 47500
$ws.Range("L74").Value = 47500
$ws.Range("N74").Value = -49248
$ws.Range("H77").Value = 47500
$ws.Range("J77").Value = 47500
$ws.Range("L77").Value = 142500
$ws.Range("N77").Value = -151236
$ws.Range("H111").Value = 50000
$ws.Range("J111").Value = 50000
$ws.Range("L111").Value = 50000
$ws.Range("N111").Value = -58180

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 661.5
$ws.Range("I109").Value = 661.5
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 1984.5
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -944.5
$ws.Range("N109").ClearContents()
$ws.Range("H132").Value = 3212.5
$ws.Range("H134").Value = 1833.3334
$ws.Range("I134").Value = 1833.3334
$ws.Range("K134").Value = 5500.0002
$ws.Range("M134").Value = -430.0002000000004
$ws.Range("H140").Value = 2499.1428
$ws.Range("I140").Value = 2499.1428
$ws.Range("K140").Value = 7497.428400000001
$ws.Range("M140").Value = -2317.428400000001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 107.666664
$ws.Range("I4").Value = 107.666664
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 107.666664
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 4.333336000000003
$ws.Range("N4").ClearContents()
$ws.Range("H14").Value = 833376
$ws.Range("I14").Value = 1035768.44
$ws.Range("J14").Value = 125002.5
$ws.Range("K14").Value = 1035768.44
$ws.Range("L14").Value = 125002.5
$ws.Range("M14").Value = -1035600.44
$ws.Range("N14").Value = -125338.5
$ws.Range("H70").Value = 5976.75
$ws.Range("I70").Value = 5976.75
$ws.Range("K70").Value = 5976.75
$ws.Range("M70").Value = -5706.75
$ws.Range("H73").Value = 5976.75
$ws.Range("I73").Value = 5976.75
$ws.Range("K73").Value = 5976.75
$ws.Range("M73").Value = -5040.75

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1755
$ws.Range("I9").Value = 1755
$ws.Range("K9").Value = 1755
$ws.Range("M9").Value = -1531

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 472
$ws.Range("I17").Value = 401.33334
$ws.Range("J17").Value = 525
$ws.Range("K17").Value = 401.33334
$ws.Range("L17").Value = 525
$ws.Range("M17").Value = -229.33334
$ws.Range("N17").Value = -869
$ws.Range("H69").Value = 18825
$ws.Range("J69").Value = 18825
$ws.Range("L69").Value = 18825
$ws.Range("N69").Value = -20323
$ws.Range("H72").Value = 18825
$ws.Range("J72").Value = 18825
$ws.Range("L72").Value = 56475
$ws.Range("N72").Value = -63963
$ws.Range("H100").Value = 300.5
$ws.Range("J100").Value = 303
$ws.Range("L100").Value = 606
$ws.Range("N100").Value = -1688
$ws.Range("H132").Value = 1830.3334
$ws.Range("I132").Value = 1245.75
$ws.Range("K132").Value = 3737.25
$ws.Range("M132").Value = -1207.25

Write-Host "Edit applied successfully"